$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (columns D = Price, E = Volume(1h))
# Values are kept as literal text (matching the source inlineStr cells),
# so each cell is pre-formatted as Text before the value is written --
# otherwise Excel would auto-convert these numeric-looking strings/percentages
# into Number cells.
$updates = [ordered]@{
    "D2" = "285.09"
    "E2" = "2.88%"
    "D3" = "28.61"
    "E3" = "5.15%"
    "D4" = "4.918"
    "E4" = "1.32%"
    "D5" = "0.06488"
    "E5" = "1.22%"
    "D6" = "7.220"
    "E6" = "4.15%"
    "D7" = "1.339"
    "E7" = "12.58%"
    "E8" = "4.32%"
    "D9" = "0.1543"
    "E9" = "0.23%"
    "D10" = "0.06364"
    "E10" = "24.06%"
    "D11" = "0.07645"
    "E11" = "2.05%"
    "E12" = "0.76%"
    "D13" = "0.08966"
    "E13" = "-0.12%"
    "D14" = "0.001599"
    "E14" = "2.38%"
    "D15" = "0.0006525"
    "E15" = "2.36%"
    "D16" = "0.006068"
    "E16" = "-0.52%"
    "D17" = "3.462"
    "E17" = "-0.42%"
    "D18" = "3.367"
    "E18" = "1.71%"
    "D19" = "2.244"
    "E19" = "-1.35%"
    "D21" = "0.1341"
    "E21" = "-0.48%"
    "D22" = "3.994"
    "E22" = "1.91%"
    "D24" = "0.04460"
    "E24" = "1.30%"
    "D25" = "0.001189"
    "E25" = "1.13%"
    "E26" = "12.08%"
    "E28" = "-9.17%"
    "E29" = "-15.73%"
    "D40" = "0.04162"
    "E40" = "0.02%"
    "D41" = "0.006707"
    "E41" = "-1.21%"
    "D42" = "0.1233"
    "E42" = "5.08%"
    "E43" = "10.91%"
    "D44" = "0.01175"
    "E44" = "2.20%"
    "D45" = "0.00005370"
    "E45" = "1.15%"
    "E46" = "-0.09%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
